# Auto-generated edit script: updates F-column "想去人数" (interested-count)
# values across all 4 sheets, per commit "Update gh-pages to output generated at 456a3b4".

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1581
$ws.Range("F3").Value = 965
$ws.Range("F4").Value = 944
$ws.Range("F5").Value = 541
$ws.Range("F6").Value = 8194
$ws.Range("F7").Value = 144
$ws.Range("F8").Value = 44
$ws.Range("F10").Value = 5888
$ws.Range("F11").Value = 588
$ws.Range("F14").Value = 8311
$ws.Range("F15").Value = 9765
$ws.Range("F16").Value = 1171
$ws.Range("F17").Value = 973
$ws.Range("F18").Value = 4640
$ws.Range("F19").Value = 722
$ws.Range("F20").Value = 308
$ws.Range("F21").Value = 91
$ws.Range("F22").Value = 301
$ws.Range("F23").Value = 166
$ws.Range("F24").Value = 1239
$ws.Range("F25").Value = 150
$ws.Range("F26").Value = 1764
$ws.Range("F27").Value = 784
$ws.Range("F28").Value = 1043
$ws.Range("F29").Value = 377
$ws.Range("F30").Value = 1945
$ws.Range("F31").Value = 362
$ws.Range("F32").Value = 514
$ws.Range("F33").Value = 2432
$ws.Range("F35").Value = 128
$ws.Range("F36").Value = 1546
$ws.Range("F37").Value = 82
$ws.Range("F39").Value = 20
$ws.Range("F40").Value = 831
$ws.Range("F41").Value = 542
$ws.Range("F43").Value = 67
$ws.Range("F44").Value = 453
$ws.Range("F45").Value = 536
$ws.Range("F48").Value = 190
$ws.Range("F49").Value = 4135

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 27

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 5513

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1581
$ws.Range("F4").Value = 965
$ws.Range("F5").Value = 944
$ws.Range("F7").Value = 8194
$ws.Range("F8").Value = 144
$ws.Range("F9").Value = 44
$ws.Range("F12").Value = 5888
$ws.Range("F13").Value = 588
$ws.Range("F14").Value = 8311
$ws.Range("F15").Value = 9765
$ws.Range("F17").Value = 1171
$ws.Range("F18").Value = 973
$ws.Range("F19").Value = 4640
$ws.Range("F20").Value = 722
$ws.Range("F21").Value = 308
$ws.Range("F22").Value = 91
$ws.Range("F23").Value = 301
$ws.Range("F25").Value = 1239
$ws.Range("F26").Value = 150
$ws.Range("F27").Value = 1764
$ws.Range("F28").Value = 784
$ws.Range("F29").Value = 1043
$ws.Range("F30").Value = 377
$ws.Range("F32").Value = 1945
$ws.Range("F33").Value = 362
$ws.Range("F34").Value = 2432
$ws.Range("F35").Value = 82
$ws.Range("F38").Value = 831
$ws.Range("F40").Value = 542
$ws.Range("F43").Value = 67
$ws.Range("F44").Value = 453
$ws.Range("F45").Value = 536
$ws.Range("F47").Value = 190
$ws.Range("F48").Value = 4135
